$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an AutoFilter on column C (the 3rd column of the A1:C51 filter range)
# so that only rows with "N" remain visible - this hides all the "Y" rows
# and marks the sheet as being in filter mode.
$null = $ws.Range("A1:C51").AutoFilter(3, @("N"), 7)

# Mark row 42 (Simon Perks / Simon Perkins) as resolved, just like the
# other previously-resolved rows (23, 40): copy the formatting from one of
# those rows so the grey "resolved" fill style is reused rather than a new
# style being created.
$ws.Range("A23:C23").Copy()
$null = $ws.Range("A42:C42").PasteSpecial(-4122)

# Add the explanatory note for why row 42 is no longer highlighted.
$ws.Range("J42").Value = "Different people - both competed simulataneously on 07/10/2008"

# Update the selection / view so the active cell is A53 with no special
# top-left scroll position.
$null = $ws.Range("A53").Select()
